# Update the "time_taken" column (F) on the existing "data" sheet with refreshed
# timestamps, then add a new "metadata" sheet describing the panel query that
# produced this export.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$timeValues = @{
    2 = "2021-10-05 14:33:44.587588"
    3 = "2021-10-05 14:33:44.587596"
    4 = "2021-10-05 14:33:44.587600"
    5 = "2021-10-05 14:33:44.587603"
    6 = "2021-10-05 14:33:44.587605"
    7 = "2021-10-05 14:33:44.587608"
    8 = "2021-10-05 14:33:44.587611"
    9 = "2021-10-05 14:33:44.587613"
    10 = "2021-10-05 14:33:44.587616"
    11 = "2021-10-05 14:33:44.587619"
    12 = "2021-10-05 14:33:44.587621"
    13 = "2021-10-05 14:33:44.587624"
    14 = "2021-10-05 14:33:44.587626"
    15 = "2021-10-05 14:33:44.587629"
    16 = "2021-10-05 14:33:44.587631"
    17 = "2021-10-05 14:33:44.587634"
    18 = "2021-10-05 14:33:44.587636"
    19 = "2021-10-05 14:33:44.587639"
    20 = "2021-10-05 14:33:44.587642"
    21 = "2021-10-05 14:33:44.587644"
    22 = "2021-10-05 14:33:44.587647"
    23 = "2021-10-05 14:33:44.587649"
    24 = "2021-10-05 14:33:44.587652"
    25 = "2021-10-05 14:33:44.587654"
    26 = "2021-10-05 14:33:44.587657"
    27 = "2021-10-05 14:33:44.587659"
    28 = "2021-10-05 14:33:44.587662"
    29 = "2021-10-05 14:33:44.587664"
    30 = "2021-10-05 14:33:44.587667"
    31 = "2021-10-05 14:33:44.587669"
    32 = "2021-10-05 14:33:44.587672"
    33 = "2021-10-05 14:33:44.587674"
    34 = "2021-10-05 14:33:44.587677"
    35 = "2021-10-05 14:33:44.587680"
    36 = "2021-10-05 14:33:44.587682"
    37 = "2021-10-05 14:33:44.587685"
    38 = "2021-10-05 14:33:44.587687"
    39 = "2021-10-05 14:33:44.587690"
    40 = "2021-10-05 14:33:44.587692"
    41 = "2021-10-05 14:33:44.587695"
    42 = "2021-10-05 14:33:44.587698"
    43 = "2021-10-05 14:33:44.587700"
    44 = "2021-10-05 14:33:44.587703"
    45 = "2021-10-05 14:33:44.587705"
    46 = "2021-10-05 14:33:44.587708"
    47 = "2021-10-05 14:33:44.587710"
    48 = "2021-10-05 14:33:44.587713"
    49 = "2021-10-05 14:33:44.587715"
    50 = "2021-10-05 14:33:44.587718"
    51 = "2021-10-05 14:33:44.587720"
    52 = "2021-10-05 14:33:44.587723"
    53 = "2021-10-05 14:33:44.587726"
    54 = "2021-10-05 14:33:44.587728"
    55 = "2021-10-05 14:33:44.587731"
    56 = "2021-10-05 14:33:44.587734"
    57 = "2021-10-05 14:33:44.587736"
    58 = "2021-10-05 14:33:44.587739"
    59 = "2021-10-05 14:33:44.587741"
    60 = "2021-10-05 14:33:44.587744"
    61 = "2021-10-05 14:33:44.587746"
    62 = "2021-10-05 14:33:44.587749"
    63 = "2021-10-05 14:33:44.587751"
    64 = "2021-10-05 14:33:44.587754"
    65 = "2021-10-05 14:33:44.587756"
    66 = "2021-10-05 14:33:44.587759"
}

foreach ($row in $timeValues.Keys) {
    $dataSheet.Cells.Item([int]$row, 6).Value = $timeValues[$row]
}

# Add the new "metadata" worksheet right after "data".
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (bold/centered/bordered header style, matching the "data" sheet's
# own header row formatting).
$headers = @{
    "B1" = "data_name"
    "C1" = "data_id"
    "D1" = "data_version"
    "E1" = "data_version_created"
    "F1" = "panel_query_time"
    "G1" = "panel_get_request"
}
foreach ($addr in $headers.Keys) {
    $metaSheet.Range($addr).Value = $headers[$addr]
}
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

# Data row.
$metaSheet.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$metaSheet.Range("B2").Value = "Disorders of immune dysregulation"
$metaSheet.Range("C2").Value = 229
$metaSheet.Range("E2").Value = "2021-08-12T10:17:50.703465Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:44.584270"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/229/?format=json"

# "0.96" must stay text (not become the number 0.96) -- force text storage via
# NumberFormat, then paste-format from a plain/unstyled cell to drop the
# number-format override again, leaving a plain default-styled text cell.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.96"
$dataSheet.Range("B2").Copy()
$metaSheet.Range("D2").PasteSpecial(-4122)

$metaSheet.Range("A1").Select()
